# Adapt to changing Ampersand compiler:
# The "Afbouwregeling inconvenienten" and "Afbouwregeling ploegendienst"
# StamItem labels used to carry a run of trailing (non-breaking) spaces
# left over from the old compiler output; the newer compiler emits them
# trimmed. Update the two cells that hold those labels accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("StamItem")

$ws.Range("B30").Value = "Afbouwregeling inconvenienten"
$ws.Range("B31").Value = "Afbouwregeling ploegendienst"

# Leave the selection where the edit was last made, like a live user edit would.
[void]$ws.Range("B31").Select()
